$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 62 (shifts rows 62..157 down to 63..158)
$ws.Rows.Item(62).Insert()

# New row 62 gets a copy of what is now row 61's content (the original row 61 data)
for ($col = 1; $col -le 20; $col++) {
    $ws.Cells.Item(62, $col).Value2 = $ws.Cells.Item(61, $col).Value2
}

# Update the date on row 61 to the new date value
$ws.Cells.Item(61, 4).Value2 = 45272
